$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$endRow = 200
$col = 3  # Column C ("Förändrad")

for ($r = $startRow; $r -le $endRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq 45202) {
        $cell.Value = 45203
    }
}
